# Update the fixed "date" placeholder text (Header & Footer -> date/time)
# from 03/12/2023 to 22/07/2024 across the slide master and every slide
# layout, matching the OOXML diff (each <a:fld type="datetimeFigureOut">
# placeholder's visible text is refreshed to the new date).

$p = $ppt.ActivePresentation
$newDate = "22/07/2024"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $isDatePlaceholder = $false
            try {
                if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePlaceholder = $true
                }
            } catch {
                $isDatePlaceholder = $false
            }

            if ($isDatePlaceholder) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout (CustomLayouts hang off the master in the object model)
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholder $layout.Shapes
}
